$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "varFile"
$ws.Range("E1").Value = "varBenar"

$ws.Range("C2").Value = 94
$ws.Range("D2").Value = "D:\Template-Upload-Rewards.xlsx"
$ws.Range("E2").Value = "Y"

$ws.Range("C3").Value = 1100
$ws.Range("D3").Value = "D:\Template-Upload-Rewards.xlsx"
$ws.Range("E3").Value = "N"

$ws.Columns.Item(1).ColumnWidth = 22.833333333333332
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 33.5

$ws.Range("E6").Select()
